$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("E2").Value = 962
$ws.Range("E3").Value = 607
$ws.Range("E4").Value = 803
$ws.Range("E5").Value = 963
$ws.Range("E6").Value = 925
$ws.Range("E7").Value = 935
$ws.Range("E8").Value = 999
